$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1096.4
$ws.Range("I19").Value = 980
$ws.Range("J19").Value = 1125.5
$ws.Range("K19").Value = 980
$ws.Range("L19").Value = 1125.5
$ws.Range("M19").Value = -805
$ws.Range("N19").Value = -1475.5
$ws.Range("H111").Value = 58824668
$ws.Range("I111").Value = 591.25
$ws.Range("J111").Value = 111112740
$ws.Range("K111").Value = 1773.75
$ws.Range("L111").Value = 333338220
$ws.Range("M111").Value = 1293.25
$ws.Range("N111").Value = -333344354
$ws.Range("H113").Value = 6863.8066
$ws.Range("I113").Value = 2435.7144
$ws.Range("K113").Value = 2435.7144
$ws.Range("M113").Value = 818.2856000000002
$ws.Range("H141").Value = 3878.6
$ws.Range("I141").Value = 1328.2162
$ws.Range("J141").Value = 35333.332
$ws.Range("K141").Value = 3984.6486
$ws.Range("L141").Value = 105999.996
$ws.Range("M141").Value = 1195.3514
$ws.Range("N141").Value = -116359.996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1123.7715
$ws.Range("I2").Value = 1191.5358
$ws.Range("J2").Value = 852.7143
$ws.Range("K2").Value = 1191.5358
$ws.Range("L2").Value = 852.7143
$ws.Range("M2").Value = -1078.5358
$ws.Range("N2").Value = -1078.7143
$ws.Range("H74").Value = 3551.7021
$ws.Range("I74").Value = 4153.5674
$ws.Range("K74").Value = 4153.5674
$ws.Range("M74").Value = -3279.5674
$ws.Range("H77").Value = 3551.7021
$ws.Range("I77").Value = 4153.5674
$ws.Range("K77").Value = 20767.837
$ws.Range("M77").Value = -16399.837
$ws.Range("H110").Value = 18936
$ws.Range("I110").Value = 21522
$ws.Range("K110").Value = 21522
$ws.Range("M110").Value = -19477
$ws.Range("H116").Value = 1123.7715
$ws.Range("I116").Value = 1191.5358
$ws.Range("J116").Value = 852.7143
$ws.Range("K116").Value = 1191.5358
$ws.Range("L116").Value = 852.7143
$ws.Range("M116").Value = 1102.4642
$ws.Range("N116").Value = -5440.7143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1123.7715
$ws.Range("I3").Value = 1191.5358
$ws.Range("J3").Value = 852.7143
$ws.Range("K3").Value = 1191.5358
$ws.Range("L3").Value = 852.7143
$ws.Range("M3").Value = -1077.5358
$ws.Range("N3").Value = -1080.7143
$ws.Range("H86").Value = 2949
$ws.Range("I86").Value = 3485.9167
$ws.Range("J86").Value = 2028.5714
$ws.Range("K86").Value = 3485.9167
$ws.Range("L86").Value = 2028.5714
$ws.Range("M86").Value = -2362.9167
$ws.Range("N86").Value = -4274.5714
$ws.Range("H89").Value = 2949
$ws.Range("I89").Value = 3485.9167
$ws.Range("J89").Value = 2028.5714
$ws.Range("K89").Value = 17429.5835
$ws.Range("L89").Value = 10142.857
$ws.Range("M89").Value = -11813.5835
$ws.Range("N89").Value = -21374.857
$ws.Range("H97").Value = 9665.6
$ws.Range("I97").Value = 5582
$ws.Range("K97").Value = 5582
$ws.Range("M97").Value = -4591
$ws.Range("H107").Value = 916.6667
$ws.Range("I107").Value = 916.6667
$ws.Range("K107").Value = 916.6667
$ws.Range("M107").Value = 1003.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1433.2157
$ws.Range("I58").Value = 1033.9487
$ws.Range("J58").Value = 2730.8333
$ws.Range("K58").Value = 1033.9487
$ws.Range("L58").Value = 2730.8333
$ws.Range("M58").Value = -830.9486999999999
$ws.Range("N58").Value = -3136.8333
$ws.Range("H99").Value = 36638
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 36638
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 36638
$ws.Range("N99").Value = -39634
$ws.Range("M99").ClearContents()
$ws.Range("H126").Value = 36638
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 36638
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 109914
$ws.Range("N126").Value = -114854
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 3929.7856
$ws.Range("I132").Value = 2801.8
$ws.Range("J132").Value = 6749.75
$ws.Range("K132").Value = 8405.400000000001
$ws.Range("L132").Value = 20249.25
$ws.Range("M132").Value = -5875.400000000001
$ws.Range("N132").Value = -25309.25
$ws.Range("H134").Value = 3155.5789
$ws.Range("I134").Value = 3340.5334
$ws.Range("K134").Value = 10021.6002
$ws.Range("M134").Value = -7486.600199999999
$ws.Range("H136").Value = 1433.2157
$ws.Range("I136").Value = 1033.9487
$ws.Range("J136").Value = 2730.8333
$ws.Range("K136").Value = 3101.8461
$ws.Range("L136").Value = 8192.499899999999
$ws.Range("M136").Value = -551.8460999999998
$ws.Range("N136").Value = -13292.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 488148
$ws.Range("I5").Value = 443.22223
$ws.Range("K5").Value = 1329.66669
$ws.Range("M5").Value = -1217.66669
$ws.Range("H18").Value = 829.8
$ws.Range("I18").Value = 605.46155
$ws.Range("J18").Value = 1246.4286
$ws.Range("K18").Value = 1816.38465
$ws.Range("L18").Value = 3739.2858
$ws.Range("M18").Value = -1647.38465
$ws.Range("N18").Value = -4077.2858
$ws.Range("H113").Value = 1528.5
$ws.Range("I113").Value = 1689.7778
$ws.Range("K113").Value = 5069.3334
$ws.Range("M113").Value = -2899.3334
$ws.Range("H135").Value = 488148
$ws.Range("I135").Value = 443.22223
$ws.Range("K135").Value = 3989.00007
$ws.Range("M135").Value = -1454.00007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2088.0232
$ws.Range("I122").Value = 1659.8846
$ws.Range("J122").Value = 2742.8235
$ws.Range("K122").Value = 4979.6538
$ws.Range("L122").Value = 8228.470499999999
$ws.Range("M122").Value = -2529.6538
$ws.Range("N122").Value = -13128.4705
$ws.Range("H126").Value = 2052.3333
$ws.Range("I126").Value = 1876.5
$ws.Range("J126").Value = 2404
$ws.Range("K126").Value = 5629.5
$ws.Range("L126").Value = 7212
$ws.Range("M126").Value = -3159.5
$ws.Range("N126").Value = -12152
$ws.Range("H132").Value = 2200.2188
$ws.Range("I132").Value = 1960.8334
$ws.Range("J132").Value = 2508
$ws.Range("K132").Value = 5882.5002
$ws.Range("L132").Value = 7524
$ws.Range("M132").Value = -3352.5002
$ws.Range("N132").Value = -12584

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7111.4365
$ws.Range("I132").Value = 7843.5405
$ws.Range("J132").Value = 5606.5557
$ws.Range("K132").Value = 23530.6215
$ws.Range("L132").Value = 16819.6671
$ws.Range("M132").Value = -21000.6215
$ws.Range("N132").Value = -21879.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1514.6086
$ws.Range("I81").Value = 1533.8462
$ws.Range("K81").Value = 3067.6924
$ws.Range("M81").Value = -2006.6924
$ws.Range("H84").Value = 1514.6086
$ws.Range("I84").Value = 1533.8462
$ws.Range("K84").Value = 15338.462
$ws.Range("M84").Value = -10034.462
$ws.Range("H132").Value = 1402.6615
$ws.Range("I132").Value = 1246.6981
$ws.Range("K132").Value = 3740.0943
$ws.Range("M132").Value = -1210.0943
$ws.Range("H136").Value = 4140
$ws.Range("I136").Value = 585.5
$ws.Range("J136").Value = 13834.091
$ws.Range("K136").Value = 1756.5
$ws.Range("L136").Value = 41502.273
$ws.Range("M136").Value = 793.5
$ws.Range("N136").Value = -46602.273
